$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Episode 15"

# Match the page margins used by the other sheets in the workbook
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

$c = $ws.Cells.Item(1, 2)
$c.Value = "Amadu"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 3)
$c.Value = "Calvin"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 4)
$c.Value = "Fabio"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 5)
$c.Value = "Luca"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 6)
$c.Value = "Lukas"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 7)
$c.Value = "Martin"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 8)
$c.Value = "Maurice"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 9)
$c.Value = "Max"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 10)
$c.Value = "Michael"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 11)
$c.Value = "Pharrell"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(1, 12)
$c.Value = "Felix"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(2, 1)
$c.Value = "Anna"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(2, 2)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(2, 3)
$c.Value = 0.06944444444444445
$c.Interior.Color = 16772625
$c.Font.Color = 0

$c = $ws.Cells.Item(2, 4)
$c.Value = 0.25
$c.Interior.Color = 16760640
$c.Font.Color = 0

$c = $ws.Cells.Item(2, 5)
$c.Value = 0.0625
$c.Interior.Color = 16772880
$c.Font.Color = 0

$c = $ws.Cells.Item(2, 6)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(2, 7)
$c.Value = 0.04166666666666666
$c.Interior.Color = 16774410
$c.Font.Color = 0

$c = $ws.Cells.Item(2, 8)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(2, 9)
$c.Value = 0.0763888888888889
$c.Interior.Color = 16772115
$c.Font.Color = 0

$c = $ws.Cells.Item(2, 10)
$c.Value = 0.4444444444444444
$c.Interior.Color = 16748145
$c.Font.Color = 15856113

$c = $ws.Cells.Item(2, 11)
$c.Value = 0.05555555555555555
$c.Interior.Color = 16773390
$c.Font.Color = 0

$c = $ws.Cells.Item(2, 12)
$c.Value = 0.09722222222222222
$c.Interior.Color = 16770840
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 1)
$c.Value = "Cecilia"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(3, 2)
$c.Value = 1
$c.Interior.Color = 16711935
$c.Font.Color = 15856113

$c = $ws.Cells.Item(3, 3)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 4)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 5)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 6)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 7)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 8)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 9)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 10)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 11)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(3, 12)
$c.Value = 0.09722222222222222
$c.Interior.Color = 16770840
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 1)
$c.Value = "Celina"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(4, 2)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 3)
$c.Value = 0.1597222222222222
$c.Interior.Color = 16766760
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 4)
$c.Value = 0.3888888888888889
$c.Interior.Color = 16751715
$c.Font.Color = 15856113

$c = $ws.Cells.Item(4, 5)
$c.Value = 0.02083333333333333
$c.Interior.Color = 16775685
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 6)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 7)
$c.Value = 0.1041666666666667
$c.Interior.Color = 16770330
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 8)
$c.Value = 0.006944444444444444
$c.Interior.Color = 16776705
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 9)
$c.Value = 0.1111111111111111
$c.Interior.Color = 16769820
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 10)
$c.Value = 0.02083333333333333
$c.Interior.Color = 16775685
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 11)
$c.Value = 0.1875
$c.Interior.Color = 16764720
$c.Font.Color = 0

$c = $ws.Cells.Item(4, 12)
$c.Value = 0.09722222222222222
$c.Interior.Color = 16770840
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 1)
$c.Value = "Franziska"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(5, 2)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 3)
$c.Value = 0.125
$c.Interior.Color = 16768800
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 4)
$c.Value = 0.08333333333333333
$c.Interior.Color = 16771605
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 5)
$c.Value = 0.09027777777777778
$c.Interior.Color = 16771095
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 6)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 7)
$c.Value = 0.2569444444444444
$c.Interior.Color = 16760385
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 8)
$c.Value = 0.01388888888888889
$c.Interior.Color = 16776195
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 9)
$c.Value = 0.1111111111111111
$c.Interior.Color = 16769820
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 10)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(5, 11)
$c.Value = 0.3194444444444444
$c.Interior.Color = 16756305
$c.Font.Color = 15856113

$c = $ws.Cells.Item(5, 12)
$c.Value = 0.09722222222222222
$c.Interior.Color = 16770840
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 1)
$c.Value = "Gina"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(6, 2)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 3)
$c.Value = 0.1458333333333333
$c.Interior.Color = 16767525
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 4)
$c.Value = 0.02777777777777778
$c.Interior.Color = 16775175
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 5)
$c.Value = 0.08333333333333333
$c.Interior.Color = 16771605
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 6)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 7)
$c.Value = 0.2291666666666667
$c.Interior.Color = 16762170
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 8)
$c.Value = 0.006944444444444444
$c.Interior.Color = 16776705
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 9)
$c.Value = 0.0625
$c.Interior.Color = 16772880
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 10)
$c.Value = 0.2222222222222222
$c.Interior.Color = 16762680
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 11)
$c.Value = 0.2222222222222222
$c.Interior.Color = 16762680
$c.Font.Color = 0

$c = $ws.Cells.Item(6, 12)
$c.Value = 0.2222222222222222
$c.Interior.Color = 16762680
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 1)
$c.Value = "Isabelle"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(7, 2)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 3)
$c.Value = 0.1458333333333333
$c.Interior.Color = 16767525
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 4)
$c.Value = 0.02777777777777778
$c.Interior.Color = 16775175
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 5)
$c.Value = 0.09027777777777778
$c.Interior.Color = 16771095
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 6)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 7)
$c.Value = 0.3541666666666667
$c.Interior.Color = 16754010
$c.Font.Color = 15856113

$c = $ws.Cells.Item(7, 8)
$c.Value = 0.05555555555555555
$c.Interior.Color = 16773390
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 9)
$c.Value = 0.1736111111111111
$c.Interior.Color = 16765740
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 10)
$c.Value = 0.1527777777777778
$c.Interior.Color = 16767015
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 11)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(7, 12)
$c.Value = 0.09722222222222222
$c.Interior.Color = 16770840
$c.Font.Color = 0

$c = $ws.Cells.Item(8, 1)
$c.Value = "Karina"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(8, 2)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(8, 3)
$c.Value = 0.2013888888888889
$c.Interior.Color = 16763955
$c.Font.Color = 0

$c = $ws.Cells.Item(8, 4)
$c.Value = 0.08333333333333333
$c.Interior.Color = 16771605
$c.Font.Color = 0

$c = $ws.Cells.Item(8, 5)
$c.Value = 0.375
$c.Interior.Color = 16752480
$c.Font.Color = 15856113

$c = $ws.Cells.Item(8, 6)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(8, 7)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(8, 8)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(8, 9)
$c.Value = 0.3263888888888889
$c.Interior.Color = 16755795
$c.Font.Color = 15856113

$c = $ws.Cells.Item(8, 10)
$c.Value = 0.01388888888888889
$c.Interior.Color = 16776195
$c.Font.Color = 0

$c = $ws.Cells.Item(8, 11)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(8, 12)
$c.Value = 0.09722222222222222
$c.Interior.Color = 16770840
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 1)
$c.Value = "Luisa"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(9, 2)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 3)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 4)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 5)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 6)
$c.Value = 1
$c.Interior.Color = 16711935
$c.Font.Color = 15856113

$c = $ws.Cells.Item(9, 7)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 8)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 9)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 10)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 11)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(9, 12)
$c.Value = 0.09722222222222222
$c.Interior.Color = 16770840
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 1)
$c.Value = "Ricarda"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(10, 2)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 3)
$c.Value = 0.01388888888888889
$c.Interior.Color = 16776195
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 4)
$c.Value = 0.05555555555555555
$c.Interior.Color = 16773390
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 5)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 6)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 7)
$c.Value = 0.01388888888888889
$c.Interior.Color = 16776195
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 8)
$c.Value = 0.9166666666666666
$c.Interior.Color = 16717290
$c.Font.Color = 15856113

$c = $ws.Cells.Item(10, 9)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 10)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 11)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(10, 12)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 1)
$c.Value = "Zoe"
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4160
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2

$c = $ws.Cells.Item(11, 2)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 3)
$c.Value = 0.1388888888888889
$c.Interior.Color = 16768035
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 4)
$c.Value = 0.08333333333333333
$c.Interior.Color = 16771605
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 5)
$c.Value = 0.2777777777777778
$c.Interior.Color = 16758855
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 6)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 7)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 8)
$c.Value = 0
$c.Interior.Color = 13882323
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 9)
$c.Value = 0.1388888888888889
$c.Interior.Color = 16768035
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 10)
$c.Value = 0.1458333333333333
$c.Interior.Color = 16767525
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 11)
$c.Value = 0.2152777777777778
$c.Interior.Color = 16762935
$c.Font.Color = 0

$c = $ws.Cells.Item(11, 12)
$c.Value = 0.09722222222222222
$c.Interior.Color = 16770840
$c.Font.Color = 0

# Keep the original first sheet as the active one (matches the target workbook state)
$wb.Worksheets.Item(1).Activate()

Write-Host "Episode 15 sheet populated"